$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 50
$ws.Range("I6").Value = 50
$ws.Range("K6").Value = 150
$ws.Range("M6").Value = -38
$ws.Range("H18").Value = 279.5
$ws.Range("I18").Value = 279.5
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 279.5
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 4.5
$ws.Range("H33").Value = 324.2
$ws.Range("I33").Value = 333.4762
$ws.Range("K33").Value = 333.4762
$ws.Range("M33").Value = -104.4762
$ws.Range("H86").Value = 2057.2104
$ws.Range("I86").Value = 1210.375
$ws.Range("J86").Value = 2673.0908
$ws.Range("K86").Value = 1210.375
$ws.Range("L86").Value = 2673.0908
$ws.Range("M86").Value = -87.375
$ws.Range("N86").Value = -4919.0908
$ws.Range("H89").Value = 2057.2104
$ws.Range("I89").Value = 1210.375
$ws.Range("J89").Value = 2673.0908
$ws.Range("K89").Value = 6051.875
$ws.Range("L89").Value = 13365.454
$ws.Range("M89").Value = -435.875
$ws.Range("N89").Value = -24597.454
$ws.Range("H96").Value = 549
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("H97").Value = 1581.4286
$ws.Range("J97").Value = 1581.4286
$ws.Range("L97").Value = 4744.2858
$ws.Range("N97").Value = -5736.2858
$ws.Range("H112").Value = 1299.2559
$ws.Range("J112").Value = 1318.2858
$ws.Range("L112").Value = 3954.8574
$ws.Range("N112").Value = -6170.857400000001
$ws.Range("H127").Value = 897.5
$ws.Range("I127").Value = 645
$ws.Range("K127").Value = 1935
$ws.Range("M127").Value = 3025
$ws.Range("H135").Value = 48388228
$ws.Range("I135").Value = 21740500
$ws.Range("J135").Value = 125000450
$ws.Range("K135").Value = 195664500
$ws.Range("L135").Value = 1125004050
$ws.Range("M135").Value = -195661965
$ws.Range("N135").Value = -1125009120
$ws.Range("H137").Value = 498203.34
$ws.Range("I137").Value = 1310.5526
$ws.Range("J137").Value = 1547199.2
$ws.Range("K137").Value = 3931.6578
$ws.Range("L137").Value = 4641597.6
$ws.Range("M137").Value = -1381.6578
$ws.Range("N137").Value = -4646697.6
$ws.Range("H138").Value = 3082.3157
$ws.Range("I138").Value = 1702.0312
$ws.Range("J138").Value = 3783.4126
$ws.Range("K138").Value = 5106.0936
$ws.Range("L138").Value = 11350.2378
$ws.Range("M138").Value = 33.90639999999985
$ws.Range("N138").Value = -21630.2378
$ws.Range("N18").ClearContents()
$ws.Range("N96").ClearContents()

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1055.1082
$ws.Range("I2").Value = 1058.5358
$ws.Range("J2").Value = 1044.4445
$ws.Range("K2").Value = 1058.5358
$ws.Range("L2").Value = 1044.4445
$ws.Range("M2").Value = -945.5358000000001
$ws.Range("N2").Value = -1270.4445
$ws.Range("H32").Value = 6845.283
$ws.Range("I32").Value = 5513.7886
$ws.Range("J32").Value = 15500
$ws.Range("K32").Value = 5513.7886
$ws.Range("L32").Value = 15500
$ws.Range("M32").Value = -5226.7886
$ws.Range("N32").Value = -16074
$ws.Range("H63").Value = 1636.909
$ws.Range("I63").Value = 1636.909
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1636.909
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -950.9090000000001
$ws.Range("H66").Value = 1636.909
$ws.Range("I66").Value = 1636.909
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 8184.545
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -4752.545
$ws.Range("H74").Value = 5761.8857
$ws.Range("I74").Value = 2929.8276
$ws.Range("J74").Value = 19450.166
$ws.Range("K74").Value = 2929.8276
$ws.Range("L74").Value = 19450.166
$ws.Range("M74").Value = -2055.8276
$ws.Range("N74").Value = -21198.166
$ws.Range("H77").Value = 5761.8857
$ws.Range("I77").Value = 2929.8276
$ws.Range("J77").Value = 19450.166
$ws.Range("K77").Value = 14649.138
$ws.Range("L77").Value = 97250.83
$ws.Range("M77").Value = -10281.138
$ws.Range("N77").Value = -105986.83
$ws.Range("H116").Value = 1055.1082
$ws.Range("I116").Value = 1058.5358
$ws.Range("J116").Value = 1044.4445
$ws.Range("K116").Value = 1058.5358
$ws.Range("L116").Value = 1044.4445
$ws.Range("M116").Value = 1235.4642
$ws.Range("N116").Value = -5632.4445
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1055.1082
$ws.Range("I3").Value = 1058.5358
$ws.Range("J3").Value = 1044.4445
$ws.Range("K3").Value = 1058.5358
$ws.Range("L3").Value = 1044.4445
$ws.Range("M3").Value = -944.5358000000001
$ws.Range("N3").Value = -1272.4445
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("H94").Value = 980.2045000000001
$ws.Range("I94").Value = 949.78125
$ws.Range("J94").Value = 1061.3334
$ws.Range("K94").Value = 949.78125
$ws.Range("L94").Value = 1061.3334
$ws.Range("M94").Value = -498.78125
$ws.Range("N94").Value = -1963.3334
$ws.Range("H99").Value = 2131.3157
$ws.Range("I99").Value = 2099.9092
$ws.Range("J99").Value = 2174.5
$ws.Range("K99").Value = 2099.9092
$ws.Range("L99").Value = 2174.5
$ws.Range("M99").Value = -601.9092000000001
$ws.Range("N99").Value = -5170.5
$ws.Range("H140").Value = 44178.06
$ws.Range("J140").Value = 44178.06
$ws.Range("L140").Value = 44178.06
$ws.Range("N140").Value = -54538.06
$ws.Range("M22").ClearContents()

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1820608.4
$ws.Range("I58").Value = 2842407
$ws.Range("J58").Value = 4077.3333
$ws.Range("K58").Value = 2842407
$ws.Range("L58").Value = 4077.3333
$ws.Range("M58").Value = -2842204
$ws.Range("N58").Value = -4483.3333
$ws.Range("H134").Value = 3118.9434
$ws.Range("I134").Value = 2057.0967
$ws.Range("J134").Value = 4615.1816
$ws.Range("K134").Value = 6171.2901
$ws.Range("L134").Value = 13845.5448
$ws.Range("M134").Value = -3636.2901
$ws.Range("N134").Value = -18915.5448
$ws.Range("H136").Value = 1820608.4
$ws.Range("I136").Value = 2842407
$ws.Range("J136").Value = 4077.3333
$ws.Range("K136").Value = 8527221
$ws.Range("L136").Value = 12231.9999
$ws.Range("M136").Value = -8524671
$ws.Range("N136").Value = -17331.9999

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 375.5
$ws.Range("I41").Value = 167.33333
$ws.Range("K41").Value = 501.99999
$ws.Range("M41").Value = -163.99999
$ws.Range("H82").Value = 3823.5334
$ws.Range("I82").Value = 1412
$ws.Range("J82").Value = 5029.3
$ws.Range("K82").Value = 4236
$ws.Range("L82").Value = 15087.9
$ws.Range("M82").Value = -3830
$ws.Range("N82").Value = -15899.9
$ws.Range("H85").Value = 3823.5334
$ws.Range("I85").Value = 1412
$ws.Range("J85").Value = 5029.3
$ws.Range("K85").Value = 4236
$ws.Range("L85").Value = 15087.9
$ws.Range("M85").Value = -2832
$ws.Range("N85").Value = -17895.9
$ws.Range("H113").Value = 719.28
$ws.Range("I113").Value = 749.6709
$ws.Range("J113").Value = 604.9524
$ws.Range("K113").Value = 2249.0127
$ws.Range("L113").Value = 1814.8572
$ws.Range("M113").Value = -79.01269999999977
$ws.Range("N113").Value = -6154.8572
$ws.Range("H132").Value = 1508.0667
$ws.Range("I132").Value = 1727.8334
$ws.Range("J132").Value = 1361.5555
$ws.Range("K132").Value = 15550.5006
$ws.Range("L132").Value = 12253.9995
$ws.Range("M132").Value = -13020.5006
$ws.Range("N132").Value = -17313.9995

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6554.5454
$ws.Range("I70").Value = 6000
$ws.Range("J70").Value = 6580.952
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 6580.952
$ws.Range("M70").Value = -5730
$ws.Range("N70").Value = -7120.952
$ws.Range("H73").Value = 6554.5454
$ws.Range("I73").Value = 6000
$ws.Range("J73").Value = 6580.952
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 6580.952
$ws.Range("M73").Value = -5064
$ws.Range("N73").Value = -8452.952000000001
$ws.Range("H80").Value = 5919.619
$ws.Range("I80").Value = 18533.334
$ws.Range("J80").Value = 3817.3333
$ws.Range("K80").Value = 18533.334
$ws.Range("L80").Value = 3817.3333
$ws.Range("M80").Value = -17535.334
$ws.Range("N80").Value = -5813.3333
$ws.Range("H83").Value = 5919.619
$ws.Range("I83").Value = 18533.334
$ws.Range("J83").Value = 3817.3333
$ws.Range("K83").Value = 92666.67
$ws.Range("L83").Value = 19086.6665
$ws.Range("M83").Value = -87674.67
$ws.Range("N83").Value = -29070.6665

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 852.8823
$ws.Range("I22").Value = 687.375
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 687.375
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -392.375
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 852.8823
$ws.Range("I27").Value = 687.375
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 687.375
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -580.375
$ws.Range("N27").Value = -1214
$ws.Range("H46").Value = 600
$ws.Range("I46").Value = 400
$ws.Range("J46").Value = 680
$ws.Range("K46").Value = 400
$ws.Range("L46").Value = 680
$ws.Range("M46").Value = -212
$ws.Range("N46").Value = -1056
$ws.Range("H55").Value = 286385
$ws.Range("I55").Value = 572071.4399999999
$ws.Range("J55").Value = 698.5714
$ws.Range("K55").Value = 572071.4399999999
$ws.Range("L55").Value = 698.5714
$ws.Range("M55").Value = -571898.4399999999
$ws.Range("N55").Value = -1044.5714
$ws.Range("H68").Value = 1268.125
$ws.Range("I68").Value = 1236.3636
$ws.Range("J68").Value = 1338
$ws.Range("K68").Value = 1236.3636
$ws.Range("L68").Value = 1338
$ws.Range("M68").Value = -487.3635999999999
$ws.Range("N68").Value = -2836
$ws.Range("H71").Value = 1268.125
$ws.Range("I71").Value = 1236.3636
$ws.Range("J71").Value = 1338
$ws.Range("K71").Value = 6181.817999999999
$ws.Range("L71").Value = 6690
$ws.Range("M71").Value = -2437.817999999999
$ws.Range("N71").Value = -14178

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1416.125
$ws.Range("I96").Value = 1434.75
$ws.Range("J96").Value = 1397.5
$ws.Range("K96").Value = 1434.75
$ws.Range("L96").Value = 1397.5
$ws.Range("M96").Value = -61.75
$ws.Range("N96").Value = -4143.5
$ws.Range("H113").Value = 1003.4211
$ws.Range("I113").Value = 448.375
$ws.Range("J113").Value = 1407.091
$ws.Range("K113").Value = 1345.125
$ws.Range("L113").Value = 4221.272999999999
$ws.Range("M113").Value = 824.875
$ws.Range("N113").Value = -8561.272999999999
$ws.Range("H136").Value = 4184.3105
$ws.Range("I136").Value = 1864.1
$ws.Range("J136").Value = 6670.25
$ws.Range("K136").Value = 5592.299999999999
$ws.Range("L136").Value = 20010.75
$ws.Range("M136").Value = -3042.299999999999
$ws.Range("N136").Value = -25110.75
$ws.Range("H139").Value = 69827.22
$ws.Range("J139").Value = 69827.22
$ws.Range("L139").Value = 69827.22
$ws.Range("N139").Value = -80107.22
